# docs(init): ajuste na matriz SWOT, 5W2H e documentos de requisitos
#
# Updates two cells in the 5W2H matrix (sheet "Plan1"):
#   F3 ("Quando" for the "Backend/scraping" row): clarify the planning
#       guidance to mention building reusable methods.
#   E6 ("Onde" for the "mobile WebView" row): narrow the target platform
#       from "Android/iOS" down to just "Android".
# Also moves the active selection to E6, matching where the author's
# cursor ended up when they saved the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "Ao planejar o sistema (fazer métodos reutilizavéis)"
$ws.Range("E6").Value = "Android"

$ws.Range("E6").Select() | Out-Null
